$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated sample data: Product name for row 2 changed from "Dryer (WM123)" to "Dryer (DR123)"
$ws.Range("B2").Value = "Dryer (DR123)"

# Updated SearchLine: move the saved cell selection from E7 to B7
$ws.Activate()
$ws.Range("B7").Select()
